# "overwrite old files with RMI version"
# The RMI copy of this workbook adds a state label ("New Mexico") next to
# the existing "Date of access" style header on the About sheet, and
# refreshes the access date shown alongside it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# B1: new label identifying the state this copy of the template is for.
$ws.Range("B1").Value = "New Mexico"

# C1: refreshed "date of access" (serial date 44509 = 2021-11-09).
$ws.Range("C1").Value = 44509
